$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray _GoBack bookmark that sits on the page-break paragraph
#    right before the "Tentative Course Schedule" table.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) "Lectures: Web API" -> "Lectures: " + "Dependency Injection" (two runs)
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("Lectures: Web API")
if ($found) {
    $para = $r.Paragraphs(1)
    $pr = $para.Range
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="439C7A99" w14:textId="77777777" w:rsidR="00A248B6" w:rsidRPr="00DC592E" w:rsidRDefault="00A248B6" w:rsidP="0069184C"><w:pPr><w:widowControl/><w:autoSpaceDE/><w:autoSpaceDN/><w:adjustRightInd/><w:rPr><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r w:rsidRPr="00DC592E"><w:rPr><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Lectures: </w:t></w:r><w:r><w:rPr><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Dependency Injection</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $pr.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 3) "Lectures: DI and Unit Testing" -> "Lectures: " + _GoBack bookmark +
#    "Unit Testing" (two runs, with the bookmark now living between them)
# ---------------------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute("Lectures: DI and Unit Testing")
if ($found2) {
    $para2 = $r2.Paragraphs(1)
    $pr2 = $para2.Range
    $xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="6CAF4DD3" w14:textId="77777777" w:rsidR="00A248B6" w:rsidRPr="00DC592E" w:rsidRDefault="00A248B6" w:rsidP="0069184C"><w:pPr><w:widowControl/><w:autoSpaceDE/><w:autoSpaceDN/><w:adjustRightInd/><w:rPr><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r w:rsidRPr="00DC592E"><w:rPr><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Lectures: </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Unit Testing</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $pr2.InsertXML($xml2)
}

# ---------------------------------------------------------------------------
# 4) Footer PAGE field cached result "2" -> "5" (document grew by a page)
# ---------------------------------------------------------------------------
$sec = $d.Sections(1)
$footer = $sec.Footers(1)
$frng = $footer.Range
$xml3 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/footer1.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.footer+xml"><pkg:xmlData><w:ftr xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:p w14:paraId="638FFE70" w14:textId="77777777" w:rsidR="002C43FE" w:rsidRDefault="002C43FE" w:rsidP="00F03EC7"><w:pPr><w:pStyle w:val="Footer"/><w:ind w:right="720"/><w:jc w:val="right"/></w:pPr><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> PAGE   \* MERGEFORMAT </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidR="0011646D"><w:rPr><w:noProof/></w:rPr><w:t>5</w:t></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r></w:p></w:ftr></pkg:xmlData></pkg:part></pkg:package>'
$frng.InsertXML($xml3)
